$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Discounted Total"
$ws.Range("B7").Value = 17.1

$ws.Columns.Item(1).ColumnWidth = 15

$ws.Range("A8").Select()
